$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark that sits after " pháp lý:" in
#    the document preamble (it will be re-created later at the point of
#    the last substantive edit, mirroring what Word itself does).
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2) First (non-bold) pricing table: "Gói 1 năm (366 ngày)" -> "Gói 12
#    tháng (366 ngày)" and "Gói 3 năm (1089 ngày)" -> "Gói 36 tháng
#    (1098 ngày)", plus the matching price bump 15.000.000 -> 13.500.000.
#    wdReplaceOne (1) is used throughout together with ranges reset to
#    the full document so each call only ever touches the single
#    occurrence intended - this table's text appears first in the
#    document, so "first match from the top" lands here.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Start = 0
$r1.End = $d.Content.End
$r1.Find.Execute("1 năm (366 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "12 tháng (366 ngày)", 1) | Out-Null

$r2 = $d.Content
$r2.Start = 0
$r2.End = $d.Content.End
$r2.Find.Execute("Gói 3 năm (1089 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "Gói 36 tháng (1098 ngày)", 1) | Out-Null

$r3 = $d.Content
$r3.Start = 0
$r3.End = $d.Content.End
$r3.Find.Execute("15.000.000 VNĐ", $false, $false, $false, $false, $false, $true, 1, $false, "13.500.000 VNĐ", 1) | Out-Null

# ---------------------------------------------------------------------
# 3) Second (bold) table further down the document repeats the same two
#    label edits (price stays as-is here). Searching $d.Content again
#    (unrestricted Start/End) continues forward past the spots already
#    rewritten above, so this now lands on the remaining/second match.
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("1 năm (366 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "12 tháng (366 ngày)", 1) | Out-Null

$r5 = $d.Content
$r5.Find.Execute("Gói 3 năm (1089 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "Gói 36 tháng (1098 ngày)", 1) | Out-Null

# ---------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark right at the point of that last
#    edit, between "1098" and " ngày)", matching where Word leaves the
#    cursor-position bookmark after the final typed change.
# ---------------------------------------------------------------------
$last = $d.Content
$last.Start = 0
$found1098 = $last.Find.Execute("Gói 36 tháng (1098 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found1098) {
    $prevEnd = $last.End
    $last.Start = $prevEnd
    $last.End = $d.Content.End
    $foundNext = $last.Find.Execute("Gói 36 tháng (1098 ngày)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $foundNext) {
        $markerPos = $prevEnd - 6
        $markerRange = $d.Range($markerPos, $markerPos)
        $d.Bookmarks.Add("_GoBack", $markerRange)
        break
    }
}

# ---------------------------------------------------------------------
# 5) Footer page-number field's cached text: "1" -> "5".
# ---------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerChars = $footer.Range.Characters
$footerChars.Item(1).Text = "5"
